$d = $word.ActiveDocument

# The document currently ends with (1-indexed paragraphs):
#   ... "What is the best design? [..sequence diagram]" paragraph
#   empty paragraph
#   paragraph containing the sequence diagram picture
#   empty paragraph                                         <- remove entirely
#   "What is the class diagram? [..class diagram]" paragraph <- remove entirely
#   paragraph containing bookmark + class diagram picture    <- keep paragraph, clear its content

$n = $d.Paragraphs.Count

# Remove the trailing bookmark/drawing content from the last paragraph,
# leaving an empty (still centered) paragraph behind.
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Range.Delete()

# Remove the whole paragraph that introduces the "What is the class diagram?" text.
$classParaIndex = $d.Paragraphs.Count
$classPara = $d.Paragraphs.Item($classParaIndex)
$classPara.Range.Delete()

# Remove the now trailing empty paragraph that used to sit between the
# sequence-diagram picture and the "What is the class diagram?" paragraph.
$emptyParaIndex = $d.Paragraphs.Count
$emptyPara = $d.Paragraphs.Item($emptyParaIndex)
$emptyPara.Range.Delete()
